$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "RM 232" (row 26). After this, all subsequent rows shift up by one,
# so what was "SC 92" (row 28) becomes row 27.
$ws.Rows.Item(26).Delete()

# Delete the row for "SC 92" (now row 27, after the previous deletion).
$ws.Rows.Item(27).Delete()

# Update the previously-blank column E (header "D") value for "SC 232" (now row 33) to -10.7.
$ws.Cells.Item(33, 5).Value = -10.7

$wb.Save()
